# Convert the three M2Doc "mergefield-style" field codes in this template
# into plain literal text using the {m:...} token syntax, as done by the
# TokenIteratorFieldRewriterSplit parser update:
#   { m:for v | null }  ->  {m:for v | null}
#   { m:v.name }        ->  {m:v.name}
#   { m:endfor }        ->  {m:endfor}
#
# Each field (fldChar begin / instrText* / fldChar end) is deleted and
# replaced by plain w:t runs containing the same text wrapped in "{" "}".
# The hidden "_GoBack" bookmark that sits inside the first field (right
# after the "null" instrText, before the closing fldChar) must keep its
# original place in the content stream.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Field 1 (paragraph 2): " " & field("m:for v | null") & bookmark _GoBack
# ---------------------------------------------------------------------
$f1 = $d.Fields.Item(1)
$f1Start = $f1.Code.Start - 1
$f1.Delete()

# Insert the replacement text plus a one-character placeholder ("X") so the
# bookmark we are about to (re)create does not land exactly on the
# paragraph-mark boundary (a position at which bookmark creation is
# unreliable). The placeholder is stripped again right after.
$ins1 = $d.Range($f1Start, $f1Start)
$ins1.Text = "{m:for v | null}X"

$bmPos = $ins1.End - 1
$d.Bookmarks.Add("_GoBack", $d.Range($bmPos, $bmPos))
$d.Range($ins1.End - 1, $ins1.End).Delete()

# ---------------------------------------------------------------------
# Field 2 (paragraph 3): "name = " & field(" m:v.name ") & ","
# ---------------------------------------------------------------------
$f2 = $d.Fields.Item(1)
$f2Start = $f2.Code.Start - 1
$f2.Delete()

$ins2 = $d.Range($f2Start, $f2Start)
$ins2.Text = "{m:v.name}"

# ---------------------------------------------------------------------
# Field 3 (paragraph 4): field(" m:endfor ")
# ---------------------------------------------------------------------
$f3 = $d.Fields.Item(1)
$f3Start = $f3.Code.Start - 1
$f3.Delete()

$ins3 = $d.Range($f3Start, $f3Start)
$ins3.Text = "{m:endfor}"
